$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-11 with refreshed weekly price data ---

# Row 2
$ws.Range("D2").Value = 44390
$ws.Range("H2").Value = 'Española'
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 16000
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = 16000
$ws.Range("N2").Value = '$/caja 30 unidades'
$ws.Range("P2").Value = 533
$ws.Range("Q2").Value = 30

# Row 3
$ws.Range("D3").Value = 44386
$ws.Range("H3").Value = 'Española'
$ws.Range("J3").Value = 30
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 15000
$ws.Range("N3").Value = '$/caja 30 unidades'
$ws.Range("P3").Value = 500
$ws.Range("Q3").Value = 30

# Row 4
$ws.Range("D4").Value = 44166
$ws.Range("H4").Value = 'Madrigal'
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = 10000
$ws.Range("N4").Value = '$/caja 40 unidades'
$ws.Range("P4").Value = 250
$ws.Range("Q4").Value = 40

# Row 5 is unchanged

# Row 6
$ws.Range("D6").Value = 44418
$ws.Range("H6").Value = 'Española'
$ws.Range("J6").Value = 80
$ws.Range("K6").Value = 16000
$ws.Range("L6").Value = 16000
$ws.Range("M6").Value = 16000
$ws.Range("N6").Value = '$/caja 30 unidades'
$ws.Range("P6").Value = 533
$ws.Range("Q6").Value = 30

# Row 7
$ws.Range("D7").Value = 44407
$ws.Range("H7").Value = 'Española'
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 18000
$ws.Range("L7").Value = 18000
$ws.Range("M7").Value = 18000
$ws.Range("N7").Value = '$/caja 30 unidades'
$ws.Range("P7").Value = 600
$ws.Range("Q7").Value = 30

# Row 8
$ws.Range("D8").Value = 44421
$ws.Range("H8").Value = 'Española'
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 16500
$ws.Range("L8").Value = 16500
$ws.Range("M8").Value = 16500
$ws.Range("N8").Value = '$/caja 30 unidades'
$ws.Range("P8").Value = 550
$ws.Range("Q8").Value = 30

# Row 9
$ws.Range("D9").Value = 44176
$ws.Range("H9").Value = 'Madrigal'
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 11000
$ws.Range("L9").Value = 11000
$ws.Range("M9").Value = 11000
$ws.Range("N9").Value = '$/caja 40 unidades'
$ws.Range("P9").Value = 275
$ws.Range("Q9").Value = 40

# Row 10
$ws.Range("D10").Value = 44446
$ws.Range("H10").Value = 'Española'
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 16000
$ws.Range("L10").Value = 16000
$ws.Range("M10").Value = 16000
$ws.Range("N10").Value = '$/caja 30 unidades'
$ws.Range("P10").Value = 400
$ws.Range("Q10").Value = 30

# Row 11
$ws.Range("D11").Value = 44161
$ws.Range("H11").Value = 'Madrigal'
$ws.Range("J11").Value = 30
$ws.Range("K11").Value = 11000
$ws.Range("L11").Value = 11000
$ws.Range("M11").Value = 11000
$ws.Range("N11").Value = '$/caja 40 unidades'
$ws.Range("P11").Value = 275
$ws.Range("Q11").Value = 40

# --- Append new row 12 ---
$ws.Range("A12").Value = 4
$ws.Range("B12").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C12").Value = 'Los Lagos'
$ws.Range("D12").Value = 44400
$ws.Range("D12").NumberFormat = $ws.Range("D11").NumberFormat
$ws.Range("E12").Value = 10
$ws.Range("F12").Value = 100112013
$ws.Range("G12").Value = 'Alcachofa'
$ws.Range("H12").Value = 'Española'
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 70
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 15000
$ws.Range("N12").Value = '$/caja 30 unidades'
$ws.Range("O12").Value = 'Provincia de Limarí'
$ws.Range("P12").Value = 500
$ws.Range("Q12").Value = 30
$ws.Range("R12").Value = 'Hortaliza'
